$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the date label in A3
$ws.Range("A3").Value = "Date:31.05.19"

# Update sales figures for Robiul (row 6), Shohel (row 7), Sodor (row 8)
$ws.Range("B6").Value = 86545
$ws.Range("C6").Value = 62
$ws.Range("D6").Value = 7

$ws.Range("B7").Value = 31820
$ws.Range("C7").Value = 18
$ws.Range("D7").Value = 5

$ws.Range("B8").Value = 94660
$ws.Range("C8").Value = 75
$ws.Range("D8").Value = 8

# Update the active cell selection
$ws.Range("D11").Select()
